# Adding test case TestCase_E38 to the "Test Cases" sheet (row 39).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate the new row's cell values. The order below matters: it controls
# the order new entries are appended to the shared-strings table so it lines
# up with the target workbook (Description, then TCID, then Jira id).
$ws.Range("C39").Value = "Verify that a user's public watchlist is not visible to another user once that particular watchlist is deleted."
$ws.Range("A39").Value = "TestCase_E38"
$ws.Range("B39").Value = "OPQA-1105"
$ws.Range("D39").Value = "Y"
$ws.Range("E39").Value = "PASS"

# Copy the cell formatting from the row above (row 37 uses the same plain
# style as the new row: columns A/B/D/E share one style, column C has the
# wrap-text style) so no new style entries are introduced.
$ws.Range("A37").Copy()
$ws.Range("A39").PasteSpecial(-4122)

$ws.Range("B37").Copy()
$ws.Range("B39").PasteSpecial(-4122)

$ws.Range("C37").Copy()
$ws.Range("C39").PasteSpecial(-4122)

$ws.Range("D37").Copy()
$ws.Range("D39").PasteSpecial(-4122)

$ws.Range("E37").Copy()
$ws.Range("E39").PasteSpecial(-4122)

# Move the active selection to the new row's last populated cell, matching
# where the editor left the cursor after adding the row.
$ws.Range("D38").Select()
